$d = $word.ActiveDocument

# 1) In the closing paragraph under "Malgrupp", the author replaced the
#    word "bullshit" with "prat om annat".
$d.Content.Find.Execute("bullshit", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "prat om annat", 2)

# 2) Append a new closing sentence about the technology used, then start a
#    new "Testning" heading paragraph, followed by the (still empty)
#    paragraph that used to hold the trailing _GoBack bookmark.
$d.Content.Find.Execute("fortsätter läsa.", $false, $false, $false, $false, $false, `
                         $true, 1, $false, `
                         "fortsätter läsa. HTML5 och CSS3 användes under detta arbete.^pTestning^p", 2)

# 3) Style the new "Testning" paragraph as a level-2 heading (Rubrik2),
#    matching the other section headings in the document.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Testning*") {
        $p.Style = "Rubrik2"
    }
}
